# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-name suffixes to "_FV2210" / "_FV2304"
# - Freeze the header row
# - Turn the data range into an Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row cells (row 1) -----------------------------------
# Columns A:J used the "_old" suffix, columns L:U used the "_new" suffix;
# column K ("diff") is unchanged.
$newHeadersFV2210 = @("Segmentname_FV2210","Segmentgruppe_FV2210","Segment_FV2210","Datenelement_FV2210","Segment ID_FV2210","Code_FV2210","Qualifier_FV2210","Beschreibung_FV2210","Bedingungsausdruck_FV2210","Bedingung_FV2210")

for ($i = 0; $i -lt $newHeadersFV2210.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeadersFV2210[$i]
}

$newHeadersFV2304 = @("Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304","Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304")

for ($i = 0; $i -lt $newHeadersFV2304.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeadersFV2304[$i]
}

# --- 2. Freeze the top row --------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the used range into an Excel Table --------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U70"), $null, 1)
$tbl.Name = "Table1"
